# Applies the portfolio.xlsx update:
#  - Portfolio sheet: A2:A10 ticker codes become text (were numbers)
#  - Portfolio sheet: new row 11 added for ticker 0053L0
#  - History sheet: Total_Asset doubled, Memo text expanded

$wb = $excel.ActiveWorkbook
$portfolio = $wb.Worksheets.Item("Portfolio")
$history   = $wb.Worksheets.Item("History")

# --- Portfolio: convert existing ticker numbers (A2:A10) to plain text ---
$tickers = @("394670", "292150", "483420", "245350", "469160", "220130", "419430", "105010", "455960")

for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $i + 2
    $cell = $portfolio.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $tickers[$i]
    $cell.Style = "Normal"
}

# --- Portfolio: append new holding in row 11 ---
$portfolio.Cells.Item(11, 1).NumberFormat = "@"
$portfolio.Cells.Item(11, 1).Value = "0053L0"
$portfolio.Cells.Item(11, 1).Style = "Normal"

$portfolio.Cells.Item(11, 2).Value = "TIGER 차이나휴머노이드로봇"
$portfolio.Cells.Item(11, 3).Value = "중국본토주식"
$portfolio.Cells.Item(11, 4).Value = 19
$portfolio.Cells.Item(11, 5).Value = 0

# --- History: update total asset and memo ---
$history.Range("B2").Value = 20000000
$history.Range("D2").Value = "Initial Capital Setup"
